# Update monthly contribution figures (columns J-O) and the per-state
# total (column Q) for rows 2-52 on Sheet1, per the "Finish Webb, Christie,
# and Bush" data refresh. The state-total column previously held a
# SUM(F:P) formula (shared across rows); it is being replaced everywhere
# with a literal rolled-up value, so setting .Value (not .Formula) clears
# the formula and leaves a plain number, matching the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rowData = @(
    @{ Row=2; J=$null; K=20; L=350; M=875; N=1275; O=1125; Q=3645 },
    @{ Row=3; J=$null; K=100; L=$null; M=200; N=$null; O=$null; Q=300 },
    @{ Row=4; J=$null; K=250; L=42200; M=8050; N=1650; O=150; Q=57850 },
    @{ Row=5; J=$null; K=2700; L=5400; M=5400; N=500; O=$null; Q=22600 },
    @{ Row=6; J=42900; K=55600; L=73625; M=167025; N=29391; O=6575; Q=434666 },
    @{ Row=7; J=3200; K=18109; L=10125; M=1175; N=8900; O=14400; Q=64539 },
    @{ Row=8; J=65529.62; K=28300; L=58100; M=20425; N=11425; O=13275; Q=233325 },
    @{ Row=9; J=250; K=100; L=1500; M=$null; N=675; O=50; Q=5775 },
    @{ Row=10; J=9600; K=750; L=3550; M=11350; N=4350; O=50; Q=45150 },
    @{ Row=11; J=29500; K=12453; L=25700; M=87350; N=13626.25; O=7300; Q=269829.25 },
    @{ Row=12; J=10303; K=2750; L=3325; M=6325; N=2775; O=4950; Q=35128 },
    @{ Row=13; J=$null; K=$null; L=$null; M=$null; N=$null; O=$null; Q=0 },
    @{ Row=14; J=$null; K=$null; L=$null; M=$null; N=2871; O=200; Q=3071 },
    @{ Row=15; J=44000; K=22913; L=1110; M=4725; N=34725; O=8650; Q=146373 },
    @{ Row=16; J=1000; K=17200; L=500; M=1525; N=4300; O=$null; Q=24525 },
    @{ Row=17; J=7425; K=1000; L=2700; M=3000; N=38850; O=5650; Q=58625 },
    @{ Row=18; J=$null; K=50; L=125; M=50; N=100; O=2750; Q=8475 },
    @{ Row=19; J=$null; K=$null; L=$null; M=5200; N=200; O=150; Q=5550 },
    @{ Row=20; J=$null; K=$null; L=$null; M=$null; N=1000; O=$null; Q=3700 },
    @{ Row=21; J=500; K=5400; L=1600; M=400; N=300; O=250; Q=8700 },
    @{ Row=22; J=31503; K=7255; L=4425; M=28100; N=60600; O=10400; Q=234683 },
    @{ Row=23; J=46570; K=21350; L=7850; M=17300; N=1875; O=15425; Q=156120 },
    @{ Row=24; J=500; K=50; L=3200; M=5325; N=5400; O=250; Q=21625 },
    @{ Row=25; J=25800; K=103; L=500; M=400; N=1250; O=3900; Q=41128 },
    @{ Row=26; J=$null; K=$null; L=$null; M=$null; N=$null; O=$null; Q=3700 },
    @{ Row=27; J=21700; K=35700; L=5000; M=1701; N=-2555; O=3815; Q=70761 },
    @{ Row=28; J=$null; K=$null; L=900; M=1700; N=50; O=$null; Q=3400 },
    @{ Row=29; J=$null; K=$null; L=$null; M=$null; N=3000; O=$null; Q=3000 },
    @{ Row=30; J=200; K=500; L=5900; M=4000; N=900; O=500; Q=15150 },
    @{ Row=31; J=6500; K=300; L=800; M=9400; N=9325; O=8095.16; Q=35920.16 },
    @{ Row=32; J=635556.13; K=248086; L=216070; M=613385; N=215953.22; O=106778.88; Q=3595934.23 },
    @{ Row=33; J=$null; K=$null; L=$null; M=$null; N=$null; O=$null; Q=0 },
    @{ Row=34; J=163311; K=38375; L=85575; M=191475; N=77356; O=52725; Q=835787 },
    @{ Row=35; J=5900; K=3150; L=2195; M=1050; N=6500; O=200; Q=19495 },
    @{ Row=36; J=$null; K=$null; L=$null; M=$null; N=$null; O=$null; Q=0 },
    @{ Row=37; J=300; K=200; L=400; M=4600; N=350; O=350; Q=6200 },
    @{ Row=38; J=$null; K=650; L=175; M=150; N=650; O=$null; Q=2125 },
    @{ Row=39; J=$null; K=$null; L=100; M=150; N=550; O=250; Q=8450 },
    @{ Row=40; J=26500; K=4475; L=67500; M=80050; N=93801; O=27990; Q=387441 },
    @{ Row=41; J=5400; K=$null; L=1000; M=500; N=775; O=$null; Q=7675 },
    @{ Row=42; J=1750; K=400; L=1200; M=4100; N=1150; O=450; Q=9050 },
    @{ Row=43; J=$null; K=$null; L=$null; M=500; N=$null; O=$null; Q=500 },
    @{ Row=44; J=9650; K=$null; L=500; M=3200; N=4325; O=2800; Q=41875 },
    @{ Row=45; J=128450; K=7725; L=7350; M=6975; N=30685; O=8275; Q=263205 },
    @{ Row=46; J=$null; K=200; L=$null; M=-1400; N=100; O=2700; Q=1600 },
    @{ Row=47; J=$null; K=500; L=$null; M=$null; N=$null; O=$null; Q=1500 },
    @{ Row=48; J=46090; K=18250; L=4849; M=35965; N=10301; O=22475; Q=187300 },
    @{ Row=49; J=2700; K=250; L=400; M=175; N=250; O=50; Q=4325 },
    @{ Row=50; J=$null; K=$null; L=500; M=$null; N=2000; O=$null; Q=2500 },
    @{ Row=51; J=1000; K=$null; L=1250; M=$null; N=5400; O=3050; Q=10700 },
    @{ Row=52; J=$null; K=2700; L=$null; M=$null; N=2700; O=$null; Q=8301 }
)

foreach ($r in $rowData) {
    if ($null -ne $r.J) { $ws.Range("J$($r.Row)").Value = $r.J }
    if ($null -ne $r.K) { $ws.Range("K$($r.Row)").Value = $r.K }
    if ($null -ne $r.L) { $ws.Range("L$($r.Row)").Value = $r.L }
    if ($null -ne $r.M) { $ws.Range("M$($r.Row)").Value = $r.M }
    if ($null -ne $r.N) { $ws.Range("N$($r.Row)").Value = $r.N }
    if ($null -ne $r.O) { $ws.Range("O$($r.Row)").Value = $r.O }
    # Column Q always becomes a literal value (formula removed) on every
    # data row, including rows whose total didn't otherwise change.
    $ws.Range("Q$($r.Row)").Value = $r.Q
}

# Row 53 ("grand total" row) drops its SUM formulas entirely; F53 (which
# carried no style) is cleared away completely, while G53:Q53 (styled
# cells) are cleared back to blank-but-styled cells.
$ws.Range("F53:Q53").ClearContents()

# The frozen-header / split-pane view introduced when the sheet was
# reopened for this edit: column A stays pinned, selection lands on R1.
$ws.Range("B1").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("R1").Select()
